# Updates the cryptos list (Price column D, Volume(1h) column E) for rows 2-51
# reflecting the latest GitHub Actions scrape, per the commit diff.
# D-column prices are forced to Text format before assignment so that
# numeric-looking strings (e.g. "377.82") are not auto-converted to numbers,
# matching the original inlineStr (text) cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.868.56"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.916.84"
$ws.Range("E3").Value = "  -3.23%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "377.82"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.50"
$ws.Range("E6").Value = "  -4.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  -2.49%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("E9").Value = "  -4.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.23"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.138"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0832"
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.387.36"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.81"
$ws.Range("E14").Value = "  -5.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.30"
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.927.29"
$ws.Range("E16").Value = "  -2.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.966"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "50.812.97"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.15"
$ws.Range("E19").Value = "  -10.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.06"
$ws.Range("E20").Value = "  -5.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("E21").Value = "  -6.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.99"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.71"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.87"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.12"
$ws.Range("E26").Value = "  +7.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.56"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("E28").Value = "  +7.27%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.165"
$ws.Range("E30").Value = "  -5.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.42"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.72"
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.51"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.78"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.94"
$ws.Range("E38").Value = "  -6.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.57"
$ws.Range("E39").Value = "  -4.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.41"
$ws.Range("E40").Value = "  -6.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.76"
$ws.Range("E42").Value = "  -6.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.13"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.08"
$ws.Range("E44").Value = "  -6.46%  "
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.271"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.993.72"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.18"
$ws.Range("E49").Value = "  -4.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0343"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.472"
$ws.Range("E51").Value = "  +8.68%  "
